# "Draft elec sector edits"
#
# The substantive change: on the ETS sheet, the shareweights for
# "petroleum" (row 11) and "natural gas peaker" (row 12) go from 0 to 1
# for every projection year, columns B (2020) through AF (2050).

$wb = $excel.ActiveWorkbook
$wsETS   = $wb.Worksheets.Item("ETS")
$wsAbout = $wb.Worksheets.Item("About")

$wsETS.Range("B11:AF12").Value = 1

# Reflect the saved selection state on the ETS sheet (cell U27 was
# selected there) without changing which tab/sheet is active - the
# workbook was left open on "About".
$wsETS.Activate()
$wsETS.Range("U27").Select()
$wsAbout.Activate()
